# "add data and legend"
# Adds a new column P to Sheet1 that concatenates the already-present
# GeoJSON fragment columns (I..O, plus the lon/lat columns A/B) into one
# GeoJSON "Feature" string per station row (rows 2-14), mirroring the
# existing shared-formula pattern used across the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 14

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Range("P$r").Formula = "=I$r&C$r&J$r&H$r&K$r&A$r&L$r&B$r&M$r&A$r&N$r&B$r&O$r"
}

# Match the page setup emitted for the sheet (A4 portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection/active cell moves onto the freshly-added column, and the old
# "topLeftCell" scroll anchor is dropped along with it.
$null = $ws.Range("P2:P14").Select()
